$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting of row 5 down into row 6, then overwrite values
$ws.Range("A5:F5").Copy()
$ws.Range("A6:F6").PasteSpecial(-4122)  # xlPasteFormats

# Add new row 6 data
$ws.Range("A6").Value = "Thursday"
$ws.Range("B6").Value = 45771
$ws.Range("C6").Value = 0.57291666666666663
$ws.Range("D6").Value = 0.64583333333333337
$ws.Range("E6").Formula = "=D6-C6"
$ws.Range("F6").Value = "Working on jump mechanic, camera movement"

# Convert E2:E5 into a shared formula block
$ws.Range("E2:E5").Formula = "=D2-C2"

# Update selection
$ws.Range("F12").Select()
